# Applies the commit "Adicionando todos os arquivos do projeto Univolei"
# to volei_base_dados.xlsx:
#  - equipes: add team 7 "Adversário"
#  - amistosos: close match 11 (row 12) + add new match 12 (row 13, "amistoso")
#  - sets: update set 34 score (row 35) + add new set 37 (row 38)
#  - rallies: add 5 new rally rows (252-256) for match 10 / set 1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "equipes": new row 8 -> team_id 7, team_name "Adversário"
# ---------------------------------------------------------------
$wsEquipes = $wb.Worksheets.Item("equipes")
$wsEquipes.Cells.Item(8, 1).Value = 7
$wsEquipes.Cells.Item(8, 2).Value = "Adversário"

# ---------------------------------------------------------------
# Sheet "amistosos": finish match in row 12, add match in row 13
# ---------------------------------------------------------------
$wsAmistosos = $wb.Worksheets.Item("amistosos")

# Row 12 (match_id 11) now CLOSED with finished/closed timestamps
$wsAmistosos.Cells.Item(12, 6).Value = "CLOSED"
$wsAmistosos.Cells.Item(12, 7).Value = "2025-09-30 12:34:30"
$wsAmistosos.Cells.Item(12, 8).Value = 1
$wsAmistosos.Cells.Item(12, 9).Value = "2025-09-30T12:34:30"

# Row 13 - new match_id 12
$wsAmistosos.Cells.Item(13, 1).Value = 12
$wsAmistosos.Cells.Item(13, 2).Value = 7
$wsAmistosos.Cells.Item(13, 3).Value = "'2025-09-30"
$wsAmistosos.Cells.Item(13, 4).Value = 0
$wsAmistosos.Cells.Item(13, 5).Value = 0
$wsAmistosos.Cells.Item(13, 6).Value = "CLOSED"
$wsAmistosos.Cells.Item(13, 7).Value = "2025-09-30 12:34:15"
$wsAmistosos.Cells.Item(13, 8).Value = 1
$wsAmistosos.Cells.Item(13, 9).Value = "2025-09-30T12:34:15"
$wsAmistosos.Cells.Item(13, 10).Value = "amistoso"

# ---------------------------------------------------------------
# Sheet "sets": update row 35 score, add new row 38
# ---------------------------------------------------------------
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Cells.Item(35, 4).Value = 7
$wsSets.Cells.Item(35, 5).Value = 5

$wsSets.Cells.Item(38, 1).Value = 37
$wsSets.Cells.Item(38, 2).Value = 12
$wsSets.Cells.Item(38, 3).Value = 1
$wsSets.Cells.Item(38, 4).Value = 0
$wsSets.Cells.Item(38, 5).Value = 0

# ---------------------------------------------------------------
# Sheet "rallies": append 5 new rally rows (252-256)
# ---------------------------------------------------------------
$wsRallies = $wb.Worksheets.Item("rallies")

$wsRallies.Cells.Item(252, 1).Value = 261
$wsRallies.Cells.Item(252, 2).Value = 10
$wsRallies.Cells.Item(252, 3).Value = 1
$wsRallies.Cells.Item(252, 4).Value = 8
$wsRallies.Cells.Item(252, 5).Value = "NOS"
$wsRallies.Cells.Item(252, 8).Value = "MEIO"
$wsRallies.Cells.Item(252, 9).Value = "PONTO"
$wsRallies.Cells.Item(252, 10).Value = "NOS"
$wsRallies.Cells.Item(252, 11).Value = 5
$wsRallies.Cells.Item(252, 12).Value = 3
$wsRallies.Cells.Item(252, 13).Value = "1  m"
$wsRallies.Cells.Item(252, 14).Value = "FRENTE"
$wsRallies.Cells.Item(252, 15).Value = "FRENTE"
$wsRallies.Cells.Item(252, 16).Value = "FRENTE"

$wsRallies.Cells.Item(253, 1).Value = 262
$wsRallies.Cells.Item(253, 2).Value = 10
$wsRallies.Cells.Item(253, 3).Value = 1
$wsRallies.Cells.Item(253, 4).Value = 9
$wsRallies.Cells.Item(253, 5).Value = "NOS"
$wsRallies.Cells.Item(253, 7).Value = 3
$wsRallies.Cells.Item(253, 8).Value = "MEIO"
$wsRallies.Cells.Item(253, 9).Value = "PONTO"
$wsRallies.Cells.Item(253, 10).Value = "NOS"
$wsRallies.Cells.Item(253, 11).Value = 6
$wsRallies.Cells.Item(253, 12).Value = 3
$wsRallies.Cells.Item(253, 13).Value = "1 3 m"
$wsRallies.Cells.Item(253, 14).Value = "FRENTE"
$wsRallies.Cells.Item(253, 15).Value = "FRENTE"
$wsRallies.Cells.Item(253, 16).Value = "FRENTE"

$wsRallies.Cells.Item(254, 1).Value = 263
$wsRallies.Cells.Item(254, 2).Value = 10
$wsRallies.Cells.Item(254, 3).Value = 1
$wsRallies.Cells.Item(254, 4).Value = 10
$wsRallies.Cells.Item(254, 5).Value = "NOS"
$wsRallies.Cells.Item(254, 7).Value = 3
$wsRallies.Cells.Item(254, 8).Value = "MEIO"
$wsRallies.Cells.Item(254, 9).Value = "ERRO"
$wsRallies.Cells.Item(254, 10).Value = "ADV"
$wsRallies.Cells.Item(254, 11).Value = 6
$wsRallies.Cells.Item(254, 12).Value = 4
$wsRallies.Cells.Item(254, 13).Value = "1 3 m e"
$wsRallies.Cells.Item(254, 14).Value = "FRENTE"
$wsRallies.Cells.Item(254, 15).Value = "FRENTE"
$wsRallies.Cells.Item(254, 16).Value = "FRENTE"

$wsRallies.Cells.Item(255, 1).Value = 264
$wsRallies.Cells.Item(255, 2).Value = 10
$wsRallies.Cells.Item(255, 3).Value = 1
$wsRallies.Cells.Item(255, 4).Value = 11
$wsRallies.Cells.Item(255, 5).Value = "NOS"
$wsRallies.Cells.Item(255, 7).Value = 3
$wsRallies.Cells.Item(255, 8).Value = "MEIO"
$wsRallies.Cells.Item(255, 9).Value = "PONTO"
$wsRallies.Cells.Item(255, 10).Value = "NOS"
$wsRallies.Cells.Item(255, 11).Value = 7
$wsRallies.Cells.Item(255, 12).Value = 4
$wsRallies.Cells.Item(255, 13).Value = "1 3 m"
$wsRallies.Cells.Item(255, 14).Value = "FRENTE"
$wsRallies.Cells.Item(255, 15).Value = "FRENTE"
$wsRallies.Cells.Item(255, 16).Value = "FRENTE"

$wsRallies.Cells.Item(256, 1).Value = 265
$wsRallies.Cells.Item(256, 2).Value = 10
$wsRallies.Cells.Item(256, 3).Value = 1
$wsRallies.Cells.Item(256, 4).Value = 12
$wsRallies.Cells.Item(256, 5).Value = "NOS"
$wsRallies.Cells.Item(256, 7).Value = 3
$wsRallies.Cells.Item(256, 8).Value = "MEIO"
$wsRallies.Cells.Item(256, 9).Value = "ERRO"
$wsRallies.Cells.Item(256, 10).Value = "ADV"
$wsRallies.Cells.Item(256, 11).Value = 7
$wsRallies.Cells.Item(256, 12).Value = 5
$wsRallies.Cells.Item(256, 13).Value = "1 3 m e"
$wsRallies.Cells.Item(256, 14).Value = "FRENTE"
$wsRallies.Cells.Item(256, 15).Value = "FRENTE"
$wsRallies.Cells.Item(256, 16).Value = "FRENTE"
